# Auto-generated edit script: updates cryptos list with refreshed
# prices / 1h volume percentages (GitHub Actions scrape refresh),
# including two coin-rank swaps (rows 20/21 and 48/49).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force text storage so numeric-looking strings (e.g. "519.82",
    # "0.998") are not auto-coerced to numbers by Excel, matching
    # the workbook's original inlineStr/text cell layout. The
    # NumberFormat/Style round-trip avoids leaving a stray explicit
    # cell style behind once the write is done.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "62.149.87"
Set-TextValue $ws.Range("E2") "  -3.29%  "
Set-TextValue $ws.Range("D3") "3.157.56"
Set-TextValue $ws.Range("E3") "  -5.91%  "
Set-TextValue $ws.Range("E4") "  -0.33%  "
Set-TextValue $ws.Range("D5") "519.82"
Set-TextValue $ws.Range("E5") "  -0.52%  "
Set-TextValue $ws.Range("D6") "169.39"
Set-TextValue $ws.Range("E6") "  -7.21%  "
Set-TextValue $ws.Range("D7") "0.587"
Set-TextValue $ws.Range("E7") "  -1.25%  "
Set-TextValue $ws.Range("E8") "  -0.03%  "
Set-TextValue $ws.Range("D9") "3.154.44"
Set-TextValue $ws.Range("E9") "  -5.74%  "
Set-TextValue $ws.Range("D10") "0.596"
Set-TextValue $ws.Range("E10") "  -2.32%  "
Set-TextValue $ws.Range("D11") "51.95"
Set-TextValue $ws.Range("E11") "  -7.95%  "
Set-TextValue $ws.Range("D12") "0.129"
Set-TextValue $ws.Range("E12") "  -0.53%  "
Set-TextValue $ws.Range("D13") "0.0000246"
Set-TextValue $ws.Range("E13") "  -0.55%  "
Set-TextValue $ws.Range("D14") "8.93"
Set-TextValue $ws.Range("E14") "  -1.28%  "
Set-TextValue $ws.Range("D15") "3.651.55"
Set-TextValue $ws.Range("E15") "  -6.13%  "
Set-TextValue $ws.Range("E16") "  -4.40%  "
Set-TextValue $ws.Range("D17") "3.146.54"
Set-TextValue $ws.Range("E17") "  -6.58%  "
Set-TextValue $ws.Range("D18") "16.99"
Set-TextValue $ws.Range("E18") "  -0.49%  "
Set-TextValue $ws.Range("D19") "61.926.92"
Set-TextValue $ws.Range("E19") "  -3.53%  "
Set-TextValue $ws.Range("B20") "Uniswap"
Set-TextValue $ws.Range("C20") "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D20") "10.85"
Set-TextValue $ws.Range("E20") "  +0.07%  "
Set-TextValue $ws.Range("B21") "Polygon"
Set-TextValue $ws.Range("C21") "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D21") "0.960"
Set-TextValue $ws.Range("E21") "  +1.19%  "
Set-TextValue $ws.Range("D22") "359.40"
Set-TextValue $ws.Range("E22") "  -2.02%  "
Set-TextValue $ws.Range("D23") "11.16"
Set-TextValue $ws.Range("E23") "  +5.87%  "
Set-TextValue $ws.Range("D24") "3.68"
Set-TextValue $ws.Range("E24") "  +0.69%  "
Set-TextValue $ws.Range("D25") "79.86"
Set-TextValue $ws.Range("E25") "  +0.02%  "
Set-TextValue $ws.Range("D26") "3.89"
Set-TextValue $ws.Range("E26") "  +5.69%  "
Set-TextValue $ws.Range("E27") "  +4.03%  "
Set-TextValue $ws.Range("D28") "2.59"
Set-TextValue $ws.Range("E28") "  -0.06%  "
Set-TextValue $ws.Range("D29") "11.13"
Set-TextValue $ws.Range("E29") "  +0.28%  "
Set-TextValue $ws.Range("D30") "8.04"
Set-TextValue $ws.Range("E30") "  -1.96%  "
Set-TextValue $ws.Range("D31") "633.55"
Set-TextValue $ws.Range("E31") "  -5.06%  "
Set-TextValue $ws.Range("D32") "27.93"
Set-TextValue $ws.Range("E32") "  -2.24%  "
Set-TextValue $ws.Range("D33") "6.34"
Set-TextValue $ws.Range("E33") "  -3.99%  "
Set-TextValue $ws.Range("D34") "11.19"
Set-TextValue $ws.Range("E34") "  +2.51%  "
Set-TextValue $ws.Range("E35") "  +1.58%  "
Set-TextValue $ws.Range("D36") "56.26"
Set-TextValue $ws.Range("E36") "  -4.81%  "
Set-TextValue $ws.Range("E37") "  -0.03%  "
Set-TextValue $ws.Range("D38") "36.54"
Set-TextValue $ws.Range("E38") "  +2.49%  "
Set-TextValue $ws.Range("D39") "0.369"
Set-TextValue $ws.Range("E39") "  +0.14%  "
Set-TextValue $ws.Range("D40") "0.997"
Set-TextValue $ws.Range("E40") "  -0.33%  "
Set-TextValue $ws.Range("D41") "0.0₃0692"
Set-TextValue $ws.Range("E41") "  +15.04%  "
Set-TextValue $ws.Range("E42") "  -1.62%  "
Set-TextValue $ws.Range("D43") "2.884.93"
Set-TextValue $ws.Range("E43") "  +4.70%  "
Set-TextValue $ws.Range("D44") "2.51"
Set-TextValue $ws.Range("E44") "  +11.60%  "
Set-TextValue $ws.Range("D45") "2.89"
Set-TextValue $ws.Range("E45") "  +13.93%  "
Set-TextValue $ws.Range("D46") "2.61"
Set-TextValue $ws.Range("E46") "  +0.80%  "
Set-TextValue $ws.Range("D47") "0.0386"
Set-TextValue $ws.Range("E47") "  +1.66%  "
Set-TextValue $ws.Range("B48") "ThetaToken"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue $ws.Range("D48") "2.53"
Set-TextValue $ws.Range("E48") "  -6.26%  "
Set-TextValue $ws.Range("B49") "ApeXProtocol"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue $ws.Range("D49") "2.92"
Set-TextValue $ws.Range("E49") "  +3.15%  "
Set-TextValue $ws.Range("E50") "  -1.06%  "
Set-TextValue $ws.Range("D51") "132.96"
Set-TextValue $ws.Range("E51") "  -1.13%  "

Write-Output "Applied 100 cell updates"
